# Grafo con informação das conexões e velocidade (Giga)
#
# Inserts a new "Fonte" column at the front of the PORTAS_CORE sheet,
# filling it with the source switch name "ICR01" for every existing
# data row, and updates the sheet's selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A; everything that used to be
# in A:E (header row + the 5 port rows + all the formatted placeholder
# rows below) shifts one column to the right, into B:F.
$ws.Columns("A").Insert() | Out-Null

# New column A needs to look like the rest of the table: copy the
# formatting from the corresponding cell in column B (now holding the
# original header/data formatting) before stamping in the new text.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A6").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Value = "Fonte"
$ws.Range("A2:A6").Value = "ICR01"

# Row 50 is a fully-styled blank "filler" row (A:E originally) that, in
# the saved workbook, keeps its original lettering and simply gains the
# new F cell rather than sliding right like the data rows above it.
$ws.Range("D50").Copy() | Out-Null
$ws.Range("A50").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Match the author's final on-screen selection.
$ws.Range("A3:A6").Select() | Out-Null
